# Apply the "many reading quiz updates, grade updates, and midterm study guide
# and solutions" changes to the Grades workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns: Homework 7 (N) and Midterm 2 (O) -----------------
$ws.Range("N2").Value = "Homework 7"
$ws.Range("O2").Value = "Midterm 2"

# Give the two new columns explicit widths (closest achievable values once
# the engine snaps to whole-pixel column widths, target ~12.8 and ~11.39).
$ws.Columns.Item(14).ColumnWidth = 12.0
$ws.Columns.Item(15).ColumnWidth = 10.5

# --- Row 5 gains a Homework 3 (H) score of 0 -------------------------------
$ws.Range("H5").Formula = "=0"

# --- New "Final Project Proposal" (M) marks of 1 for most students --------
$ws.Range("M3").Value = 1
$ws.Range("M4").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("M6").Value = 1
$ws.Range("M7").Value = 1
$ws.Range("M8").Value = 1
$ws.Range("M9").Value = 1
$ws.Range("M12").Value = 1
$ws.Range("M16").Value = 1
$ws.Range("M17").Value = 1
$ws.Range("M19").Value = 1

# --- Update the active cell selection to match the new edit location ------
[void]$ws.Range("O3").Select()
